# Auto-generated Excel COM-interop script
# Applies market-price data refresh values to 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per the scheduled-runner commit that updated Kujata_Profits.xlsx leve-profit rows.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 2845.3333
$ws.Range("I74").Value = 2845.3333
$ws.Range("K74").Value = 2845.3333
$ws.Range("M74").Value = -1909.3333
# Row 77
$ws.Range("H77").Value = 2845.3333
$ws.Range("I77").Value = 2845.3333
$ws.Range("K77").Value = 14226.6665
$ws.Range("M77").Value = -9546.666499999999
# Row 111
$ws.Range("H111").Value = 3915.0625
$ws.Range("I111").Value = 2163.1428
$ws.Range("J111").Value = 5277.6665
$ws.Range("K111").Value = 6489.428400000001
$ws.Range("L111").Value = 15832.9995
$ws.Range("M111").Value = -3422.428400000001
$ws.Range("N111").Value = -21966.9995

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 31
$ws.Range("H31").Value = 11621.375
$ws.Range("I31").Value = 4328.5
$ws.Range("J31").Value = 33500
$ws.Range("K31").Value = 4328.5
$ws.Range("L31").Value = 33500
$ws.Range("M31").Value = -4034.5
$ws.Range("N31").Value = -34088
# Row 32
$ws.Range("H32").Value = 3009.309
$ws.Range("I32").Value = 3420.439
$ws.Range("J32").Value = 1805.2858
$ws.Range("K32").Value = 3420.439
$ws.Range("L32").Value = 1805.2858
$ws.Range("M32").Value = -3133.439
$ws.Range("N32").Value = -2379.2858
# Row 74
$ws.Range("H74").Value = 1028.92
$ws.Range("I74").Value = 879.75
$ws.Range("K74").Value = 879.75
$ws.Range("M74").Value = -5.75
# Row 77
$ws.Range("H77").Value = 1028.92
$ws.Range("I77").Value = 879.75
$ws.Range("K77").Value = 4398.75
$ws.Range("M77").Value = -30.75
# Row 108
$ws.Range("H108").Value = 40000
$ws.Range("J108").Value = 40000
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 53
$ws.Range("H53").Value = 34878
$ws.Range("J53").Value = 34878
$ws.Range("L53").Value = 34878
$ws.Range("N53").Value = -36026
# Row 102
$ws.Range("H102").Value = 20571.428
$ws.Range("I102").Value = 12500
$ws.Range("J102").Value = 31333.334
$ws.Range("K102").Value = 12500
$ws.Range("L102").Value = 31333.334
$ws.Range("M102").Value = -9255
$ws.Range("N102").Value = -37823.334
# Row 108
$ws.Range("H108").Value = 19842
$ws.Range("J108").Value = 19842
$ws.Range("L108").Value = 19842
$ws.Range("N108").Value = -27522

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1462.3914
$ws.Range("I31").Value = 1443.25
$ws.Range("J31").Value = 1472.6
$ws.Range("K31").Value = 1443.25
$ws.Range("L31").Value = 1472.6
$ws.Range("M31").Value = -1148.25
$ws.Range("N31").Value = -2062.6
# Row 34
$ws.Range("H34").Value = 1462.3914
$ws.Range("I34").Value = 1443.25
$ws.Range("J34").Value = 1472.6
$ws.Range("K34").Value = 1443.25
$ws.Range("L34").Value = 1472.6
$ws.Range("M34").Value = -1241.25
$ws.Range("N34").Value = -1876.6
# Row 62
$ws.Range("H62").Value = 7695318.5
$ws.Range("I62").Value = 3055.652
$ws.Range("J62").Value = 66669336
$ws.Range("K62").Value = 3055.652
$ws.Range("L62").Value = 66669336
$ws.Range("M62").Value = -2431.652
$ws.Range("N62").Value = -66670584
# Row 65
$ws.Range("H65").Value = 7695318.5
$ws.Range("I65").Value = 3055.652
$ws.Range("J65").Value = 66669336
$ws.Range("K65").Value = 15278.26
$ws.Range("L65").Value = 333346680
$ws.Range("M65").Value = -12158.26
$ws.Range("N65").Value = -333352920
# Row 134
$ws.Range("H134").Value = 15874520
$ws.Range("I134").Value = 17545206
$ws.Range("K134").Value = 52635618
$ws.Range("M134").Value = -52633083
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 151633.5
$ws.Range("J4").Value = 167857.12
$ws.Range("L4").Value = 503571.36
$ws.Range("N4").Value = -503795.36
# Row 68
$ws.Range("H68").Value = 1659.1708
$ws.Range("I68").Value = 627.1
$ws.Range("J68").Value = 1992.0968
$ws.Range("K68").Value = 1881.3
$ws.Range("L68").Value = 5976.2904
$ws.Range("M68").Value = -1070.3
$ws.Range("N68").Value = -7598.2904
# Row 71
$ws.Range("H71").Value = 1659.1708
$ws.Range("I71").Value = 627.1
$ws.Range("J71").Value = 1992.0968
$ws.Range("K71").Value = 5643.900000000001
$ws.Range("L71").Value = 17928.8712
$ws.Range("M71").Value = -1587.900000000001
$ws.Range("N71").Value = -26040.8712
# Row 94
$ws.Range("H94").Value = 5140.875
$ws.Range("J94").Value = 5140.875
$ws.Range("L94").Value = 15422.625
$ws.Range("N94").Value = -16774.625
# Row 100
$ws.Range("H100").Value = 2988.5715
$ws.Range("J100").Value = 2988.5715
$ws.Range("L100").Value = 8965.7145
$ws.Range("N100").Value = -10587.7145
# Row 107
$ws.Range("H107").Value = 5306.2607
$ws.Range("J107").Value = 11050
$ws.Range("L107").Value = 33150
$ws.Range("N107").Value = -36990
# Row 112
$ws.Range("H112").Value = 6000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 6000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 18000
$ws.Range("M112").Value = $null
$ws.Range("N112").Value = -20216
# Row 121
$ws.Range("H121").Value = 281.5
$ws.Range("I121").Value = 281.5
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 844.5
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = 465.5
$ws.Range("N121").Value = $null
# Row 122
$ws.Range("H122").Value = 586.8
$ws.Range("I122").Value = 510.66666
$ws.Range("J122").Value = 701
$ws.Range("K122").Value = 4595.99994
$ws.Range("L122").Value = 6309
$ws.Range("M122").Value = -2145.99994
$ws.Range("N122").Value = -11209
# Row 123
$ws.Range("H123").Value = 100
$ws.Range("I123").Value = 100
$ws.Range("K123").Value = 300
$ws.Range("M123").Value = 2150
# Row 133
$ws.Range("H133").Value = 3407.3333
$ws.Range("I133").Value = 1109.8572
$ws.Range("J133").Value = 4211.45
$ws.Range("K133").Value = 3329.5716
$ws.Range("L133").Value = 12634.35
$ws.Range("M133").Value = 1730.4284
$ws.Range("N133").Value = -22754.35
# Row 136
$ws.Range("H136").Value = 1808.1333
$ws.Range("I136").Value = 868.4286
$ws.Range("J136").Value = 2630.375
$ws.Range("K136").Value = 2605.2858
$ws.Range("L136").Value = 7891.125
$ws.Range("M136").Value = 2494.7142
$ws.Range("N136").Value = -18091.125
# Row 137
$ws.Range("H137").Value = 6805.1514
$ws.Range("I137").Value = 1308.3846
$ws.Range("J137").Value = 10378.05
$ws.Range("K137").Value = 3925.1538
$ws.Range("L137").Value = 31134.15
$ws.Range("M137").Value = 1174.8462
$ws.Range("N137").Value = -41334.14999999999
# Row 138
$ws.Range("H138").Value = 2398.2307
$ws.Range("I138").Value = 2354.4
$ws.Range("K138").Value = 7063.200000000001
$ws.Range("M138").Value = -1923.200000000001
# Row 139
$ws.Range("H139").Value = 1759.7576
$ws.Range("I139").Value = 1615.8334
$ws.Range("J139").Value = 1932.4667
$ws.Range("K139").Value = 4847.5002
$ws.Range("L139").Value = 5797.4001
$ws.Range("M139").Value = 292.4997999999996
$ws.Range("N139").Value = -16077.4001
# Row 140
$ws.Range("H140").Value = 28798.605
$ws.Range("J140").Value = 2869.1667
$ws.Range("L140").Value = 8607.500100000001
$ws.Range("N140").Value = -18967.5001

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 1001.4286
$ws.Range("I31").Value = 1001.4286
$ws.Range("K31").Value = 1001.4286
$ws.Range("M31").Value = -709.4286
# Row 37
$ws.Range("H37").Value = 1001.4286
$ws.Range("I37").Value = 1001.4286
$ws.Range("K37").Value = 1001.4286
$ws.Range("M37").Value = -724.4286

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1803.5883
$ws.Range("I7").Value = 1655.0769
$ws.Range("J7").Value = 2286.25
$ws.Range("K7").Value = 1655.0769
$ws.Range("L7").Value = 2286.25
$ws.Range("M7").Value = -1543.0769
$ws.Range("N7").Value = -2510.25
# Row 22
$ws.Range("H22").Value = 1358.1666
$ws.Range("J22").Value = 1885.5714
$ws.Range("L22").Value = 1885.5714
$ws.Range("N22").Value = -2475.5714
# Row 27
$ws.Range("H27").Value = 1358.1666
$ws.Range("J27").Value = 1885.5714
$ws.Range("L27").Value = 1885.5714
$ws.Range("N27").Value = -2099.5714
# Row 46
$ws.Range("H46").Value = 5402.4707
$ws.Range("I46").Value = 530.25
$ws.Range("J46").Value = 9733.333000000001
$ws.Range("K46").Value = 530.25
$ws.Range("L46").Value = 9733.333000000001
$ws.Range("M46").Value = -342.25
$ws.Range("N46").Value = -10109.333
# Row 126
$ws.Range("H126").Value = 1803.5883
$ws.Range("I126").Value = 1655.0769
$ws.Range("J126").Value = 2286.25
$ws.Range("K126").Value = 4965.2307
$ws.Range("L126").Value = 6858.75
$ws.Range("M126").Value = -2495.2307
$ws.Range("N126").Value = -11798.75
# Row 136
$ws.Range("H136").Value = 6905.706
$ws.Range("I136").Value = 9199.75
$ws.Range("J136").Value = 1400
$ws.Range("K136").Value = 27599.25
$ws.Range("L136").Value = 4200
$ws.Range("M136").Value = -25049.25
$ws.Range("N136").Value = -9300

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 58
$ws.Range("H58").Value = 11698
$ws.Range("I58").Value = 10000
$ws.Range("J58").Value = 12547
$ws.Range("K58").Value = 10000
$ws.Range("L58").Value = 12547
$ws.Range("M58").Value = -9692
$ws.Range("N58").Value = -13163
# Row 123
$ws.Range("H123").Value = 58502.8
$ws.Range("J123").Value = 58502.8
$ws.Range("L123").Value = 58502.8
$ws.Range("N123").Value = -68302.8
# Row 126
$ws.Range("H126").Value = 48309790
$ws.Range("J126").Value = 405
$ws.Range("L126").Value = 1215
$ws.Range("N126").Value = -6155
